# Update the "Generate Report for Handback" timestamps.
# These cells hold their timestamp as plain text (string) values even
# though they are styled with a date-like number format, and the
# workbook preserves that text representation as long as we simply
# assign the new string values without touching NumberFormat.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 00be81ab...md (row 3)
$overview.Range("G3").Value = "2016-08-25 12:47:42"

# zh-cn sheet, row 3 (00be81ab...zh-cn.xlf): Correspond Handoff / Handback Datetime
$zhcn.Range("H3").Value = "2016-08-25 12:47:38"
$zhcn.Range("K3").Value = "2016-08-25 12:48:09"

# de-de sheet, row 3 (00be81ab...de-de.xlf): Correspond Handoff / Handback Datetime
$dede.Range("H3").Value = "2016-08-25 12:47:42"
$dede.Range("K3").Value = "2016-08-25 12:48:18"
